$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Get-ParagraphRange($doc, $searchText) {
    $findRng = $doc.Range($doc.Content.Start, $doc.Content.End)
    $found = $findRng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { throw "text not found: $searchText" }
    $para = $findRng.Paragraphs(1)
    return $doc.Range($para.Range.Start, $para.Range.End - 1)
}

# --- Change 1: "called indx.html " -> "called " / "book" / "ind" / "e" / "x.html " ---
$p1 = Get-ParagraphRange $d "called indx.html"

$body1 = '<w:p>' +
  '<w:r><w:t xml:space="preserve">Create a file </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">called </w:t></w:r>' +
  '<w:r><w:t>book</w:t></w:r>' +
  '<w:r><w:t>ind</w:t></w:r>' +
  '<w:r><w:t>e</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">x.html </w:t></w:r>' +
  '<w:r><w:t>within your IDE</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> and store it within the Vagrant shared directory under the HTML directory</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
  '</w:p>'
$p1.InsertXML($pkgOpen + $body1 + $pkgClose)

# --- Change 2: add proofErr spellStart/spellEnd around "VisualStudio" ---
$p2 = Get-ParagraphRange $d "new folder icon next to the project directory within VisualStudio code"

$body2 = '<w:p>' +
  '<w:r><w:t xml:space="preserve">To do this use the new folder icon next to the project directory within </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>VisualStudio</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> code. This is the second icon along:</w:t></w:r>' +
  '</w:p>'
$p2.InsertXML($pkgOpen + $body2 + $pkgClose)

# --- Change 3: add proofErr gramStart/gramEnd around "top level" ---
$p3 = Get-ParagraphRange $d "Create the top level body element"

$body3 = '<w:p>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Create the </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>top level</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> body element</w:t></w:r>' +
  '</w:p>'
$p3.InsertXML($pkgOpen + $body3 + $pkgClose)

Write-Output "done"
